$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.NumberFormat = "General"
    $cell.Style = $origStyle
}

Set-TextValue $ws.Cells.Item(2, 4) "42.928.28"
Set-TextValue $ws.Cells.Item(2, 5) "  +2.40%  "
Set-TextValue $ws.Cells.Item(3, 4) "2.308.64"
Set-TextValue $ws.Cells.Item(3, 5) "  +1.19%  "
Set-TextValue $ws.Cells.Item(4, 5) "  +0.03%  "
Set-TextValue $ws.Cells.Item(5, 4) "321.97"
Set-TextValue $ws.Cells.Item(5, 5) "  +2.11%  "
Set-TextValue $ws.Cells.Item(6, 4) "104.89"
Set-TextValue $ws.Cells.Item(6, 5) "  +2.35%  "
Set-TextValue $ws.Cells.Item(7, 5) "  +1.21%  "
Set-TextValue $ws.Cells.Item(8, 5) "  +0.09%  "
Set-TextValue $ws.Cells.Item(9, 5) "  +1.58%  "
Set-TextValue $ws.Cells.Item(10, 4) "40.50"
Set-TextValue $ws.Cells.Item(10, 5) "  +5.17%  "
Set-TextValue $ws.Cells.Item(11, 4) "0.0911"
Set-TextValue $ws.Cells.Item(11, 5) "  +0.93%  "
Set-TextValue $ws.Cells.Item(12, 4) "8.64"
Set-TextValue $ws.Cells.Item(12, 5) "  +5.51%  "
Set-TextValue $ws.Cells.Item(13, 5) "  +1.03%  "
Set-TextValue $ws.Cells.Item(14, 4) "0.977"
Set-TextValue $ws.Cells.Item(14, 5) "  +2.49%  "
Set-TextValue $ws.Cells.Item(15, 4) "15.41"
Set-TextValue $ws.Cells.Item(15, 5) "  +1.32%  "
Set-TextValue $ws.Cells.Item(16, 4) "2.658.49"
Set-TextValue $ws.Cells.Item(16, 5) "  +1.38%  "
Set-TextValue $ws.Cells.Item(17, 4) "2.307.53"
Set-TextValue $ws.Cells.Item(17, 5) "  +1.29%  "
Set-TextValue $ws.Cells.Item(18, 4) "42.870.14"
Set-TextValue $ws.Cells.Item(18, 5) "  +2.47%  "
Set-TextValue $ws.Cells.Item(19, 4) "7.59"
Set-TextValue $ws.Cells.Item(19, 5) "  +2.22%  "
Set-TextValue $ws.Cells.Item(20, 5) "  +1.40%  "
Set-TextValue $ws.Cells.Item(21, 4) "13.33"
Set-TextValue $ws.Cells.Item(21, 5) "  +32.87%  "
Set-TextValue $ws.Cells.Item(22, 5) "  +1.06%  "
Set-TextValue $ws.Cells.Item(23, 4) "3.60"
Set-TextValue $ws.Cells.Item(23, 5) "  +1.61%  "
Set-TextValue $ws.Cells.Item(24, 4) "272.54"
Set-TextValue $ws.Cells.Item(24, 5) "  -1.92%  "
Set-TextValue $ws.Cells.Item(25, 4) "2.25"
Set-TextValue $ws.Cells.Item(25, 5) "  +0.52%  "
Set-TextValue $ws.Cells.Item(26, 5) "  -0.35%  "
Set-TextValue $ws.Cells.Item(27, 4) "10.99"
Set-TextValue $ws.Cells.Item(27, 5) "  +2.84%  "
Set-TextValue $ws.Cells.Item(28, 5) "  +0.28%  "
Set-TextValue $ws.Cells.Item(29, 4) "22.76"
Set-TextValue $ws.Cells.Item(29, 5) "  -0.48%  "
Set-TextValue $ws.Cells.Item(30, 4) "37.99"
Set-TextValue $ws.Cells.Item(30, 5) "  +9.78%  "
Set-TextValue $ws.Cells.Item(31, 4) "165.97"
Set-TextValue $ws.Cells.Item(31, 5) "  +2.02%  "
Set-TextValue $ws.Cells.Item(32, 4) "6.20"
Set-TextValue $ws.Cells.Item(32, 5) "  +6.77%  "
Set-TextValue $ws.Cells.Item(33, 5) "  +3.02%  "
Set-TextValue $ws.Cells.Item(34, 5) "  +1.35%  "
Set-TextValue $ws.Cells.Item(35, 4) "0.116"
Set-TextValue $ws.Cells.Item(35, 5) "  +1.83%  "
Set-TextValue $ws.Cells.Item(36, 5) "  -11.37%  "
Set-TextValue $ws.Cells.Item(37, 5) "  +2.97%  "
Set-TextValue $ws.Cells.Item(38, 5) "  +3.22%  "
Set-TextValue $ws.Cells.Item(39, 5) "  +3.14%  "
Set-TextValue $ws.Cells.Item(40, 4) "2.76"
Set-TextValue $ws.Cells.Item(40, 5) "  -4.10%  "
Set-TextValue $ws.Cells.Item(41, 4) "1.58"
Set-TextValue $ws.Cells.Item(41, 5) "  +8.65%  "
Set-TextValue $ws.Cells.Item(42, 4) "102.06"
Set-TextValue $ws.Cells.Item(42, 5) "  +2.05%  "
Set-TextValue $ws.Cells.Item(43, 4) "71.14"
Set-TextValue $ws.Cells.Item(43, 5) "  +3.30%  "
Set-TextValue $ws.Cells.Item(44, 4) "0.227"
Set-TextValue $ws.Cells.Item(44, 5) "  +1.59%  "
Set-TextValue $ws.Cells.Item(45, 5) "  -0.03%  "
Set-TextValue $ws.Cells.Item(46, 4) "12.46"
Set-TextValue $ws.Cells.Item(46, 5) "  +5.82%  "
Set-TextValue $ws.Cells.Item(47, 4) "82.65"
Set-TextValue $ws.Cells.Item(47, 5) "  +9.60%  "
Set-TextValue $ws.Cells.Item(48, 4) "115.13"
Set-TextValue $ws.Cells.Item(48, 5) "  -0.27%  "
Set-TextValue $ws.Cells.Item(49, 5) "  +1.65%  "
Set-TextValue $ws.Cells.Item(50, 4) "8.90"
Set-TextValue $ws.Cells.Item(50, 5) "  -0.54%  "
Set-TextValue $ws.Cells.Item(51, 4) "1.589.44"
Set-TextValue $ws.Cells.Item(51, 5) "  +4.86%  "
